$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the existing "Yht" total row (currently row 35) down to row 37,
# making room for two new entries. Copying+clearing (rather than a
# sheet-wide Insert/shift) keeps this scoped to the B:D data columns.
$ws.Range("B35:D35").Copy($ws.Range("B37"))
$ws.Range("B35:D35").ClearContents()

# Give the two freed-up rows the same look (borders/fonts/number format
# and wrap text) as the preceding data row.
$ws.Range("B34:D34").Copy()
$ws.Range("B35:D36").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Match the row height used by the wrapped-description rows elsewhere.
$ws.Rows.Item(35).RowHeight = 37.5
$ws.Rows.Item(36).RowHeight = 37.5

# Row 35: 2024-03-07 (serial 45358), 2 hours
$ws.Cells.Item(35, 2).Value = (Get-Date -Year 2024 -Month 3 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(35, 3).Value = 2
$ws.Cells.Item(35, 4).Value = "Koitin saada asetuksiin liityvät normi sivu toiminnot toimimaan modulaarisesti"

# Row 36: 2024-03-08 (serial 45359), 5 hours
$ws.Cells.Item(36, 2).Value = (Get-Date -Year 2024 -Month 3 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(36, 3).Value = 5
$ws.Cells.Item(36, 4).Value = "Sain asetukset toimimaan oikein modulaarisesti normi sivu toiminnon osalta"

# Fix up the SUM formula in the total row, now on row 37, to cover the
# newly added rows.
$ws.Range("C37").Formula = "=SUM(C6:C36)"

# Update the selection to match the recorded post-edit state.
$ws.Range("I34").Select()

$wb.Save()
